$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z5").Value = "hi"
$ws.Range("Z5").Font.Name = "Times New Roman"
$ws.Range("Z5").Font.Family = 1
$ws.Range("Z5").Font.Size = 10
$ws.Range("Z5").Font.Bold = $true
